# Applies the scheduled-runner market/profit data refresh across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets (currentAveragePrice* / LevePrice* / LeveProfit*
# columns H:N). Values come from the latest market-board pull; a few rows gain or
# lose their LeveProfit (N or M) cell depending on whether that leve currently has an
# HQ/NQ price recorded.
$wb = $excel.ActiveWorkbook
$updatedCells = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value2 = 306.33334  # ALC!H80: 321.92856 -> 306.33334
$updatedCells++
$ws.Cells.Item(80, 9).Value2 = 314.18182  # ALC!I80: 300.83334 -> 314.18182
$updatedCells++
$ws.Cells.Item(80, 10).Value2 = 284.75  # ALC!J80: 448.5 -> 284.75
$updatedCells++
$ws.Cells.Item(80, 11).Value2 = 942.54546  # ALC!K80: 902.5000200000001 -> 942.54546
$updatedCells++
$ws.Cells.Item(80, 12).Value2 = 854.25  # ALC!L80: 1345.5 -> 854.25
$updatedCells++
$ws.Cells.Item(80, 13).Value2 = 55.45453999999995  # ALC!M80: 95.49997999999994 -> 55.45453999999995
$updatedCells++
$ws.Cells.Item(80, 14).Value2 = -2850.25  # ALC!N80: -3341.5 -> -2850.25
$updatedCells++
$ws.Cells.Item(81, 8).Value2 = 80000  # ALC!H81: 0 -> 80000
$updatedCells++
$ws.Cells.Item(81, 10).Value2 = 80000  # ALC!J81: 0 -> 80000
$updatedCells++
$ws.Cells.Item(81, 12).Value2 = 80000  # ALC!L81: 0 -> 80000
$updatedCells++
$ws.Cells.Item(81, 14).Value2 = -81996  # ALC!N81: (new) -> -81996
$updatedCells++
$ws.Cells.Item(83, 8).Value2 = 306.33334  # ALC!H83: 321.92856 -> 306.33334
$updatedCells++
$ws.Cells.Item(83, 9).Value2 = 314.18182  # ALC!I83: 300.83334 -> 314.18182
$updatedCells++
$ws.Cells.Item(83, 10).Value2 = 284.75  # ALC!J83: 448.5 -> 284.75
$updatedCells++
$ws.Cells.Item(83, 11).Value2 = 2827.63638  # ALC!K83: 2707.50006 -> 2827.63638
$updatedCells++
$ws.Cells.Item(83, 12).Value2 = 2562.75  # ALC!L83: 4036.5 -> 2562.75
$updatedCells++
$ws.Cells.Item(83, 13).Value2 = 2164.36362  # ALC!M83: 2284.49994 -> 2164.36362
$updatedCells++
$ws.Cells.Item(83, 14).Value2 = -12546.75  # ALC!N83: -14020.5 -> -12546.75
$updatedCells++
$ws.Cells.Item(84, 8).Value2 = 80000  # ALC!H84: 0 -> 80000
$updatedCells++
$ws.Cells.Item(84, 10).Value2 = 80000  # ALC!J84: 0 -> 80000
$updatedCells++
$ws.Cells.Item(84, 12).Value2 = 240000  # ALC!L84: 0 -> 240000
$updatedCells++
$ws.Cells.Item(84, 14).Value2 = -249984  # ALC!N84: (new) -> -249984
$updatedCells++
$ws.Cells.Item(113, 8).Value2 = 28583278  # ALC!H113: 26326940 -> 28583278
$updatedCells++
$ws.Cells.Item(113, 9).Value2 = 33336658  # ALC!I113: 30306324 -> 33336658
$updatedCells++
$ws.Cells.Item(113, 11).Value2 = 33336658  # ALC!K113: 30306324 -> 33336658
$updatedCells++
$ws.Cells.Item(113, 13).Value2 = -33333404  # ALC!M113: -30303070 -> -33333404
$updatedCells++
$ws.Cells.Item(116, 8).Value2 = 3164.5  # ALC!H116: 3155.4211 -> 3164.5
$updatedCells++
$ws.Cells.Item(116, 9).Value2 = 2984  # ALC!I116: 3042 -> 2984
$updatedCells++
$ws.Cells.Item(116, 10).Value2 = 3299.875  # ALC!J116: 3349.8572 -> 3299.875
$updatedCells++
$ws.Cells.Item(116, 11).Value2 = 2984  # ALC!K116: 3042 -> 2984
$updatedCells++
$ws.Cells.Item(116, 12).Value2 = 3299.875  # ALC!L116: 3349.8572 -> 3299.875
$updatedCells++
$ws.Cells.Item(116, 13).Value2 = 458  # ALC!M116: 400 -> 458
$updatedCells++
$ws.Cells.Item(116, 14).Value2 = -10183.875  # ALC!N116: -10233.8572 -> -10183.875
$updatedCells++
$ws.Cells.Item(127, 8).Value2 = 855.8333  # ALC!H127: 809.25 -> 855.8333
$updatedCells++
$ws.Cells.Item(127, 9).Value2 = 855.8333  # ALC!I127: 829.2857 -> 855.8333
$updatedCells++
$ws.Cells.Item(127, 10).Value2 = 0  # ALC!J127: 669 -> 0
$updatedCells++
$ws.Cells.Item(127, 11).Value2 = 2567.4999  # ALC!K127: 2487.8571 -> 2567.4999
$updatedCells++
$ws.Cells.Item(127, 12).Value2 = 0  # ALC!L127: 2007 -> 0
$updatedCells++
$ws.Cells.Item(127, 13).Value2 = 2392.5001  # ALC!M127: 2472.1429 -> 2392.5001
$updatedCells++
$ws.Cells.Item(127, 14).ClearContents()  # ALC!N127: -11927 -> (removed)
$updatedCells++
$ws.Cells.Item(138, 8).Value2 = 4556.6895  # ALC!H138: 4389.8066 -> 4556.6895
$updatedCells++
$ws.Cells.Item(138, 10).Value2 = 3985.25  # ALC!J138: 3802.0454 -> 3985.25
$updatedCells++
$ws.Cells.Item(138, 12).Value2 = 11955.75  # ALC!L138: 11406.1362 -> 11955.75
$updatedCells++
$ws.Cells.Item(138, 14).Value2 = -22235.75  # ALC!N138: -21686.1362 -> -22235.75
$updatedCells++

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 6049.2583  # ARM!H32: 6701.3228 -> 6049.2583
$updatedCells++
$ws.Cells.Item(32, 9).Value2 = 6049.2583  # ARM!I32: 6701.3228 -> 6049.2583
$updatedCells++
$ws.Cells.Item(32, 11).Value2 = 6049.2583  # ARM!K32: 6701.3228 -> 6049.2583
$updatedCells++
$ws.Cells.Item(32, 13).Value2 = -5762.2583  # ARM!M32: -6414.3228 -> -5762.2583
$updatedCells++
$ws.Cells.Item(44, 8).Value2 = 31000  # ARM!H44: 32000 -> 31000
$updatedCells++
$ws.Cells.Item(44, 10).Value2 = 31000  # ARM!J44: 32000 -> 31000
$updatedCells++
$ws.Cells.Item(44, 12).Value2 = 31000  # ARM!L44: 32000 -> 31000
$updatedCells++
$ws.Cells.Item(44, 14).Value2 = -31976  # ARM!N44: -32976 -> -31976
$updatedCells++
$ws.Cells.Item(61, 8).Value2 = 6899.24  # ARM!H61: 7159.125 -> 6899.24
$updatedCells++
$ws.Cells.Item(61, 9).Value2 = 6899.24  # ARM!I61: 7159.125 -> 6899.24
$updatedCells++
$ws.Cells.Item(61, 11).Value2 = 6899.24  # ARM!K61: 7159.125 -> 6899.24
$updatedCells++
$ws.Cells.Item(61, 13).Value2 = -6687.24  # ARM!M61: -6947.125 -> -6687.24
$updatedCells++
$ws.Cells.Item(63, 8).Value2 = 1115853.8  # ARM!H63: 1255010.5 -> 1115853.8
$updatedCells++
$ws.Cells.Item(63, 9).Value2 = 4415  # ARM!I63: 4778 -> 4415
$updatedCells++
$ws.Cells.Item(63, 11).Value2 = 4415  # ARM!K63: 4778 -> 4415
$updatedCells++
$ws.Cells.Item(63, 13).Value2 = -3729  # ARM!M63: -4092 -> -3729
$updatedCells++
$ws.Cells.Item(66, 8).Value2 = 1115853.8  # ARM!H66: 1255010.5 -> 1115853.8
$updatedCells++
$ws.Cells.Item(66, 9).Value2 = 4415  # ARM!I66: 4778 -> 4415
$updatedCells++
$ws.Cells.Item(66, 11).Value2 = 22075  # ARM!K66: 23890 -> 22075
$updatedCells++
$ws.Cells.Item(66, 13).Value2 = -18643  # ARM!M66: -20458 -> -18643
$updatedCells++
$ws.Cells.Item(136, 8).Value2 = 6899.24  # ARM!H136: 7159.125 -> 6899.24
$updatedCells++
$ws.Cells.Item(136, 9).Value2 = 6899.24  # ARM!I136: 7159.125 -> 6899.24
$updatedCells++
$ws.Cells.Item(136, 11).Value2 = 20697.72  # ARM!K136: 21477.375 -> 20697.72
$updatedCells++
$ws.Cells.Item(136, 13).Value2 = -18147.72  # ARM!M136: -18927.375 -> -18147.72
$updatedCells++

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value2 = 1207.6471  # BSM!H86: 1318.4706 -> 1207.6471
$updatedCells++
$ws.Cells.Item(86, 9).Value2 = 1169.4667  # BSM!I86: 1246.1428 -> 1169.4667
$updatedCells++
$ws.Cells.Item(86, 10).Value2 = 1494  # BSM!J86: 1656 -> 1494
$updatedCells++
$ws.Cells.Item(86, 11).Value2 = 1169.4667  # BSM!K86: 1246.1428 -> 1169.4667
$updatedCells++
$ws.Cells.Item(86, 12).Value2 = 1494  # BSM!L86: 1656 -> 1494
$updatedCells++
$ws.Cells.Item(86, 13).Value2 = -46.46669999999995  # BSM!M86: -123.1428000000001 -> -46.46669999999995
$updatedCells++
$ws.Cells.Item(86, 14).Value2 = -3740  # BSM!N86: -3902 -> -3740
$updatedCells++
$ws.Cells.Item(89, 8).Value2 = 1207.6471  # BSM!H89: 1318.4706 -> 1207.6471
$updatedCells++
$ws.Cells.Item(89, 9).Value2 = 1169.4667  # BSM!I89: 1246.1428 -> 1169.4667
$updatedCells++
$ws.Cells.Item(89, 10).Value2 = 1494  # BSM!J89: 1656 -> 1494
$updatedCells++
$ws.Cells.Item(89, 11).Value2 = 5847.3335  # BSM!K89: 6230.714 -> 5847.3335
$updatedCells++
$ws.Cells.Item(89, 12).Value2 = 7470  # BSM!L89: 8280 -> 7470
$updatedCells++
$ws.Cells.Item(89, 13).Value2 = -231.3334999999997  # BSM!M89: -614.7139999999999 -> -231.3334999999997
$updatedCells++
$ws.Cells.Item(89, 14).Value2 = -18702  # BSM!N89: -19512 -> -18702
$updatedCells++
$ws.Cells.Item(134, 8).Value2 = 3824.4849  # BSM!H134: 3923.5938 -> 3824.4849
$updatedCells++
$ws.Cells.Item(134, 9).Value2 = 3207.6333  # BSM!I134: 3295.724 -> 3207.6333
$updatedCells++
$ws.Cells.Item(134, 11).Value2 = 9622.8999  # BSM!K134: 9887.172 -> 9622.8999
$updatedCells++
$ws.Cells.Item(134, 13).Value2 = -7087.8999  # BSM!M134: -7352.172 -> -7087.8999
$updatedCells++

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 38465612  # CRP!H31: 40004160 -> 38465612
$updatedCells++
$ws.Cells.Item(31, 9).Value2 = 47621580  # CRP!I31: 52634316 -> 47621580
$updatedCells++
$ws.Cells.Item(31, 10).Value2 = 10550  # CRP!J31: 8666.666999999999 -> 10550
$updatedCells++
$ws.Cells.Item(31, 11).Value2 = 47621580  # CRP!K31: 52634316 -> 47621580
$updatedCells++
$ws.Cells.Item(31, 12).Value2 = 10550  # CRP!L31: 8666.666999999999 -> 10550
$updatedCells++
$ws.Cells.Item(31, 13).Value2 = -47621285  # CRP!M31: -52634021 -> -47621285
$updatedCells++
$ws.Cells.Item(31, 14).Value2 = -11140  # CRP!N31: -9256.666999999999 -> -11140
$updatedCells++
$ws.Cells.Item(34, 8).Value2 = 38465612  # CRP!H34: 40004160 -> 38465612
$updatedCells++
$ws.Cells.Item(34, 9).Value2 = 47621580  # CRP!I34: 52634316 -> 47621580
$updatedCells++
$ws.Cells.Item(34, 10).Value2 = 10550  # CRP!J34: 8666.666999999999 -> 10550
$updatedCells++
$ws.Cells.Item(34, 11).Value2 = 47621580  # CRP!K34: 52634316 -> 47621580
$updatedCells++
$ws.Cells.Item(34, 12).Value2 = 10550  # CRP!L34: 8666.666999999999 -> 10550
$updatedCells++
$ws.Cells.Item(34, 13).Value2 = -47621378  # CRP!M34: -52634114 -> -47621378
$updatedCells++
$ws.Cells.Item(34, 14).Value2 = -10954  # CRP!N34: -9070.666999999999 -> -10954
$updatedCells++
$ws.Cells.Item(58, 8).Value2 = 12472.667  # CRP!H58: 10568.223 -> 12472.667
$updatedCells++
$ws.Cells.Item(58, 9).Value2 = 12445  # CRP!I58: 5605.6 -> 12445
$updatedCells++
$ws.Cells.Item(58, 11).Value2 = 12445  # CRP!K58: 5605.6 -> 12445
$updatedCells++
$ws.Cells.Item(58, 13).Value2 = -12242  # CRP!M58: -5402.6 -> -12242
$updatedCells++
$ws.Cells.Item(94, 8).Value2 = 2023.1538  # CRP!H94: 2023.2307 -> 2023.1538
$updatedCells++
$ws.Cells.Item(94, 9).Value2 = 1656.8  # CRP!I94: 1824 -> 1656.8
$updatedCells++
$ws.Cells.Item(94, 10).Value2 = 2252.125  # CRP!J94: 2111.7778 -> 2252.125
$updatedCells++
$ws.Cells.Item(94, 11).Value2 = 1656.8  # CRP!K94: 1824 -> 1656.8
$updatedCells++
$ws.Cells.Item(94, 12).Value2 = 2252.125  # CRP!L94: 2111.7778 -> 2252.125
$updatedCells++
$ws.Cells.Item(94, 13).Value2 = -1205.8  # CRP!M94: -1373 -> -1205.8
$updatedCells++
$ws.Cells.Item(94, 14).Value2 = -3154.125  # CRP!N94: -3013.7778 -> -3154.125
$updatedCells++
$ws.Cells.Item(110, 8).Value2 = 79520  # CRP!H110: 82209 -> 79520
$updatedCells++
$ws.Cells.Item(110, 10).Value2 = 79520  # CRP!J110: 82209 -> 79520
$updatedCells++
$ws.Cells.Item(110, 12).Value2 = 79520  # CRP!L110: 82209 -> 79520
$updatedCells++
$ws.Cells.Item(110, 14).Value2 = -87700  # CRP!N110: -90389 -> -87700
$updatedCells++
$ws.Cells.Item(134, 8).Value2 = 2846.4614  # CRP!H134: 2472.8386 -> 2846.4614
$updatedCells++
$ws.Cells.Item(134, 9).Value2 = 2049.4211  # CRP!I134: 1732.875 -> 2049.4211
$updatedCells++
$ws.Cells.Item(134, 11).Value2 = 6148.263300000001  # CRP!K134: 5198.625 -> 6148.263300000001
$updatedCells++
$ws.Cells.Item(134, 13).Value2 = -3613.263300000001  # CRP!M134: -2663.625 -> -3613.263300000001
$updatedCells++
$ws.Cells.Item(136, 8).Value2 = 12472.667  # CRP!H136: 10568.223 -> 12472.667
$updatedCells++
$ws.Cells.Item(136, 9).Value2 = 12445  # CRP!I136: 5605.6 -> 12445
$updatedCells++
$ws.Cells.Item(136, 11).Value2 = 37335  # CRP!K136: 16816.8 -> 37335
$updatedCells++
$ws.Cells.Item(136, 13).Value2 = -34785  # CRP!M136: -14266.8 -> -34785
$updatedCells++

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value2 = 569.1429000000001  # CUL!H34: 749.5 -> 569.1429000000001
$updatedCells++
$ws.Cells.Item(34, 10).Value2 = 496.8  # CUL!J34: 749 -> 496.8
$updatedCells++
$ws.Cells.Item(34, 12).Value2 = 1490.4  # CUL!L34: 2247 -> 1490.4
$updatedCells++
$ws.Cells.Item(34, 14).Value2 = -1658.4  # CUL!N34: -2415 -> -1658.4
$updatedCells++
$ws.Cells.Item(39, 8).Value2 = 1625.8572  # CUL!H39: 1687.6666 -> 1625.8572
$updatedCells++
$ws.Cells.Item(39, 9).Value2 = 1729.5  # CUL!I39: 1837.8 -> 1729.5
$updatedCells++
$ws.Cells.Item(39, 10).Value2 = 1004  # CUL!J39: 1500 -> 1004
$updatedCells++
$ws.Cells.Item(39, 11).Value2 = 5188.5  # CUL!K39: 5513.4 -> 5188.5
$updatedCells++
$ws.Cells.Item(39, 12).Value2 = 3012  # CUL!L39: 4500 -> 3012
$updatedCells++
$ws.Cells.Item(39, 13).Value2 = -4894.5  # CUL!M39: -5219.4 -> -4894.5
$updatedCells++
$ws.Cells.Item(39, 14).Value2 = -3600  # CUL!N39: -5088 -> -3600
$updatedCells++
$ws.Cells.Item(49, 8).Value2 = 525  # CUL!H49: 1035 -> 525
$updatedCells++
$ws.Cells.Item(49, 9).Value2 = 800  # CUL!I49: 1558.3334 -> 800
$updatedCells++
$ws.Cells.Item(49, 11).Value2 = 2400  # CUL!K49: 4675.0002 -> 2400
$updatedCells++
$ws.Cells.Item(49, 13).Value2 = -2244  # CUL!M49: -4519.0002 -> -2244
$updatedCells++
$ws.Cells.Item(55, 8).Value2 = 7679.143  # CUL!H55: 7712.25 -> 7679.143
$updatedCells++
$ws.Cells.Item(55, 9).Value2 = 1004  # CUL!I55: 686 -> 1004
$updatedCells++
$ws.Cells.Item(55, 10).Value2 = 8791.666999999999  # CUL!J55: 17549 -> 8791.666999999999
$updatedCells++
$ws.Cells.Item(55, 11).Value2 = 3012  # CUL!K55: 2058 -> 3012
$updatedCells++
$ws.Cells.Item(55, 12).Value2 = 26375.001  # CUL!L55: 52647 -> 26375.001
$updatedCells++
$ws.Cells.Item(55, 13).Value2 = -2835  # CUL!M55: -1881 -> -2835
$updatedCells++
$ws.Cells.Item(55, 14).Value2 = -26729.001  # CUL!N55: -53001 -> -26729.001
$updatedCells++
$ws.Cells.Item(129, 8).Value2 = 654063  # CUL!H129: 653991.1 -> 654063
$updatedCells++
$ws.Cells.Item(129, 9).Value2 = 169762.17  # CUL!I129: 203270.2 -> 169762.17
$updatedCells++
$ws.Cells.Item(129, 10).Value2 = 918227.0600000001  # CUL!J129: 841791.5 -> 918227.0600000001
$updatedCells++
$ws.Cells.Item(129, 11).Value2 = 509286.51  # CUL!K129: 609810.6000000001 -> 509286.51
$updatedCells++
$ws.Cells.Item(129, 12).Value2 = 2754681.18  # CUL!L129: 2525374.5 -> 2754681.18
$updatedCells++
$ws.Cells.Item(129, 13).Value2 = -504286.51  # CUL!M129: -604810.6000000001 -> -504286.51
$updatedCells++
$ws.Cells.Item(129, 14).Value2 = -2764681.18  # CUL!N129: -2535374.5 -> -2764681.18
$updatedCells++

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value2 = 1397.7646  # GSM!H132: 1400.3529 -> 1397.7646
$updatedCells++
$ws.Cells.Item(132, 9).Value2 = 1397.7646  # GSM!I132: 1400.3529 -> 1397.7646
$updatedCells++
$ws.Cells.Item(132, 11).Value2 = 4193.293799999999  # GSM!K132: 4201.0587 -> 4193.293799999999
$updatedCells++
$ws.Cells.Item(132, 13).Value2 = -1663.293799999999  # GSM!M132: -1671.0587 -> -1663.293799999999
$updatedCells++

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value2 = 4277.6665  # LTW!H40: 4360.778 -> 4277.6665
$updatedCells++
$ws.Cells.Item(40, 9).Value2 = 3499.0527  # LTW!I40: 3617.158 -> 3499.0527
$updatedCells++
$ws.Cells.Item(40, 11).Value2 = 3499.0527  # LTW!K40: 3617.158 -> 3499.0527
$updatedCells++
$ws.Cells.Item(40, 13).Value2 = -3363.0527  # LTW!M40: -3481.158 -> -3363.0527
$updatedCells++
$ws.Cells.Item(55, 8).Value2 = 278.44446  # LTW!H55: 273.27777 -> 278.44446
$updatedCells++
$ws.Cells.Item(55, 9).Value2 = 227.90909  # LTW!I55: 237.9 -> 227.90909
$updatedCells++
$ws.Cells.Item(55, 10).Value2 = 357.85715  # LTW!J55: 317.5 -> 357.85715
$updatedCells++
$ws.Cells.Item(55, 11).Value2 = 227.90909  # LTW!K55: 237.9 -> 227.90909
$updatedCells++
$ws.Cells.Item(55, 12).Value2 = 357.85715  # LTW!L55: 317.5 -> 357.85715
$updatedCells++
$ws.Cells.Item(55, 13).Value2 = -54.90908999999999  # LTW!M55: -64.90000000000001 -> -54.90908999999999
$updatedCells++
$ws.Cells.Item(55, 14).Value2 = -703.85715  # LTW!N55: -663.5 -> -703.85715
$updatedCells++
$ws.Cells.Item(99, 8).Value2 = 20659.666  # LTW!H99: 22997.5 -> 20659.666
$updatedCells++
$ws.Cells.Item(99, 9).Value2 = 20659.666  # LTW!I99: 20663.334 -> 20659.666
$updatedCells++
$ws.Cells.Item(99, 10).Value2 = 0  # LTW!J99: 30000 -> 0
$updatedCells++
$ws.Cells.Item(99, 11).Value2 = 20659.666  # LTW!K99: 20663.334 -> 20659.666
$updatedCells++
$ws.Cells.Item(99, 12).Value2 = 0  # LTW!L99: 30000 -> 0
$updatedCells++
$ws.Cells.Item(99, 13).Value2 = -17664.666  # LTW!M99: -17668.334 -> -17664.666
$updatedCells++
$ws.Cells.Item(99, 14).ClearContents()  # LTW!N99: -35990 -> (removed)
$updatedCells++

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value2 = 10476  # WVR!H3: 10481.333 -> 10476
$updatedCells++
$ws.Cells.Item(3, 10).Value2 = 714.5  # WVR!J3: 722.5 -> 714.5
$updatedCells++
$ws.Cells.Item(3, 12).Value2 = 714.5  # WVR!L3: 722.5 -> 714.5
$updatedCells++
$ws.Cells.Item(3, 14).Value2 = -942.5  # WVR!N3: -950.5 -> -942.5
$updatedCells++
$ws.Cells.Item(81, 8).Value2 = 4191.1665  # WVR!H81: 4472.077 -> 4191.1665
$updatedCells++
$ws.Cells.Item(81, 10).Value2 = 7574.5  # WVR!J81: 7628.2 -> 7574.5
$updatedCells++
$ws.Cells.Item(81, 12).Value2 = 15149  # WVR!L81: 15256.4 -> 15149
$updatedCells++
$ws.Cells.Item(81, 14).Value2 = -17271  # WVR!N81: -17378.4 -> -17271
$updatedCells++
$ws.Cells.Item(84, 8).Value2 = 4191.1665  # WVR!H84: 4472.077 -> 4191.1665
$updatedCells++
$ws.Cells.Item(84, 10).Value2 = 7574.5  # WVR!J84: 7628.2 -> 7574.5
$updatedCells++
$ws.Cells.Item(84, 12).Value2 = 75745  # WVR!L84: 76282 -> 75745
$updatedCells++
$ws.Cells.Item(84, 14).Value2 = -86353  # WVR!N84: -86890 -> -86353
$updatedCells++
$ws.Cells.Item(96, 8).Value2 = 4837  # WVR!H96: 3044.9 -> 4837
$updatedCells++
$ws.Cells.Item(96, 9).Value2 = 4389  # WVR!I96: 2944.5 -> 4389
$updatedCells++
$ws.Cells.Item(96, 10).Value2 = 5285  # WVR!J96: 3195.5 -> 5285
$updatedCells++
$ws.Cells.Item(96, 11).Value2 = 4389  # WVR!K96: 2944.5 -> 4389
$updatedCells++
$ws.Cells.Item(96, 12).Value2 = 5285  # WVR!L96: 3195.5 -> 5285
$updatedCells++
$ws.Cells.Item(96, 13).Value2 = -3016  # WVR!M96: -1571.5 -> -3016
$updatedCells++
$ws.Cells.Item(96, 14).Value2 = -8031  # WVR!N96: -5941.5 -> -8031
$updatedCells++
$ws.Cells.Item(132, 8).Value2 = 5059.375  # WVR!H132: 4952.12 -> 5059.375
$updatedCells++
$ws.Cells.Item(132, 9).Value2 = 4627.174  # WVR!I132: 4533.4585 -> 4627.174
$updatedCells++
$ws.Cells.Item(132, 11).Value2 = 13881.522  # WVR!K132: 13600.3755 -> 13881.522
$updatedCells++
$ws.Cells.Item(132, 13).Value2 = -11351.522  # WVR!M132: -11070.3755 -> -11351.522
$updatedCells++
$ws.Cells.Item(136, 8).Value2 = 5841.263  # WVR!H136: 5622.6665 -> 5841.263
$updatedCells++
$ws.Cells.Item(136, 9).Value2 = 3611.7778  # WVR!I136: 3599.818 -> 3611.7778
$updatedCells++
$ws.Cells.Item(136, 11).Value2 = 10835.3334  # WVR!K136: 10799.454 -> 10835.3334
$updatedCells++
$ws.Cells.Item(136, 13).Value2 = -8285.3334  # WVR!M136: -8249.454000000002 -> -8285.3334
$updatedCells++

Write-Output "Updated $updatedCells cells across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
